# Apply the batch_test BOQ edits: rows 8-12 get new item data (point counts /
# descriptions shift up the table), and the Grand-Total / Net-Payable rows
# (14 & 16) pick up the newly computed amounts.
#
# Numeric-looking text (quantity/amount columns D & G that must stay text,
# e.g. "2", "25600.00") is written with the cell pre-formatted as Text ("@")
# so Excel doesn't auto-promote it to a real number - same as typing it into
# a Text-formatted cell in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "P. point" / Short point (up to 3 mtr.) ---
$ws.Range("A8").Value = "P. point"
$ws.Range("C8").Value = 100
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2"
$ws.Range("E8").Value = "Short point (up to 3 mtr.)"
$ws.Range("F8").Value = 256
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "25600.00"

# --- Row 9: "P. point" / Medium point (up to 6 mtr.) ---
$ws.Range("A9").Value = "P. point"
$ws.Range("C9").Value = 21
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3"
$ws.Range("E9").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F9").Value = 472
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "9912.00"

# --- Row 10: "Each" / flush type switch item ---
$ws.Range("A10").Value = "Each"
$ws.Range("C10").Value = 16
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.0"
$ws.Range("E10").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F10").Value = 23
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "368.00"

# --- Row 11: blank unit / "Total" row (shifted up from row 10's old role) ---
$ws.Range("A11").Value = ""
$ws.Range("C11").Value = 14
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "8"
$ws.Range("E11").Value = "Total"

# --- Row 12: "%" unit / "Add Tender Premium " row ---
$ws.Range("A12").Value = "%"
$ws.Range("C12").Value = 14
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9"
$ws.Range("E12").Value = "Add Tender Premium "

# --- Row 14: Grand Total Rs. amounts ---
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "35880.00"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "35880.00"

# --- Row 16: NET PAYABLE AMOUNT Rs. amounts ---
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "35880.00"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "35880.00"
